$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Peer  and self assessment")

# Fill in Robert's peer assessment entries (first table, row 5; second table, row 16)
$ws.Range("B5").Value = "Good"
$ws.Range("C5").Value = "Very active in our online meetings, Uploaded relevant articles to github"

$ws.Range("B16").Value = "Good"
$ws.Range("C16").Value = "Quick response in Discord, well formulated and active in meetings "
